$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 342; Excel shifts existing rows 342-448 down to 343-449
# and the worksheet dimension grows to A1:R449 automatically.
$ws.Rows(342).Insert()

# Populate the freshly-inserted row 342 with the new weekly price record.
$ws.Range("A342").Value = 5
$ws.Range("B342").Value = "Macroferia Regional de Talca"
$ws.Range("C342").Value = "Maule"
$ws.Range("D342").Value = 45093
$ws.Range("E342").Value = 7
$ws.Range("F342").Value = 100112009
$ws.Range("G342").Value = "Acelga"
$ws.Range("H342").Value = "Sin especificar"
$ws.Range("I342").Value = "Primera"
$ws.Range("J342").Value = 500
$ws.Range("K342").Value = 2000
$ws.Range("L342").Value = 2000
$ws.Range("M342").Value = 2000
$ws.Range("N342").Value = "$/docena de atados (4 kilos)"
$ws.Range("O342").Value = "Región del Maule"
$ws.Range("P342").Value = 500
$ws.Range("Q342").Value = 4
$ws.Range("R342").Value = "Hortaliza"
